$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (Sheet1 -> Libros)
$ws.Name = "Libros"

# Re-order / relabel the header row:
#   A1: ID           -> TITULO
#   B1: AUTOR        -> AUTOR   (unchanged)
#   C1: NOMBRE LIBRO -> ID
#   D1: UDS DISPONBLES -> UDS
$ws.Range("A1").Value = "TITULO"
$ws.Range("B1").Value = "AUTOR"
$ws.Range("C1").Value = "ID"
$ws.Range("D1").Value = "UDS"

# Give every header cell the same look (centered, no wrap) so the whole
# row shares one consistent style, like A1 already had.
$header = $ws.Range("A1:D1")
$header.HorizontalAlignment = -4108   # xlCenter
$header.WrapText = $false

# Leave the header row selected, matching the reviewed state of the sheet.
$header.Select()
